# Update "想去人数" (F column) counts across all sheets to reflect the
# latest generated output (commit: "Update gh-pages to output generated at 456a3b4").

$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$updates1 = @{
    2  = 119
    3  = 188
    4  = 436
    5  = 206
    6  = 137
    7  = 1210
    8  = 405
    9  = 203
    12 = 382
    13 = 422
    14 = 798
    15 = 185
    16 = 732
    17 = 293
    18 = 84
    19 = 1024
    20 = 476
    21 = 275
    25 = 47
}
foreach ($row in $updates1.Keys) {
    $ws1.Range("F$row").Value = $updates1[$row]
}

# Sheet "演出" (sheet2)
$ws2 = $wb.Worksheets.Item("演出")
$updates2 = @{
    4  = 370
    5  = 43
    8  = 85
    10 = 631
}
foreach ($row in $updates2.Keys) {
    $ws2.Range("F$row").Value = $updates2[$row]
}

# Sheet "本地生活" (sheet3)
$ws3 = $wb.Worksheets.Item("本地生活")
$updates3 = @{
    2 = 347
}
foreach ($row in $updates3.Keys) {
    $ws3.Range("F$row").Value = $updates3[$row]
}

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$updates4 = @{
    2  = 347
    4  = 119
    5  = 188
    6  = 436
    7  = 206
    8  = 137
    9  = 1210
    10 = 405
    11 = 203
    14 = 370
    16 = 43
    17 = 382
    20 = 422
    21 = 798
    22 = 185
    23 = 732
    24 = 293
    25 = 84
    26 = 1024
    27 = 476
    28 = 85
    30 = 275
    33 = 631
    36 = 47
}
foreach ($row in $updates4.Keys) {
    $ws4.Range("F$row").Value = $updates4[$row]
}
